$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1620.4
$ws.Range("I40").Value = 1550
$ws.Range("K40").Value = 1550
$ws.Range("M40").Value = -1375
$ws.Range("H64").Value = 622220.25
$ws.Range("I64").Value = 1172397.5
$ws.Range("J64").Value = 3270.875
$ws.Range("K64").Value = 1172397.5
$ws.Range("L64").Value = 3270.875
$ws.Range("M64").Value = -1172149.5
$ws.Range("N64").Value = -3766.875
$ws.Range("H67").Value = 622220.25
$ws.Range("I67").Value = 1172397.5
$ws.Range("J67").Value = 3270.875
$ws.Range("K67").Value = 1172397.5
$ws.Range("L67").Value = 3270.875
$ws.Range("M67").Value = -1171539.5
$ws.Range("N67").Value = -4986.875
$ws.Range("H76").Value = 2648962.5
$ws.Range("I76").Value = 3370374.2
$ws.Range("K76").Value = 3370374.2
$ws.Range("M76").Value = -3370059.2
$ws.Range("H79").Value = 2648962.5
$ws.Range("I79").Value = 3370374.2
$ws.Range("K79").Value = 3370374.2
$ws.Range("M79").Value = -3369282.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 3437.2222
$ws.Range("I63").Value = 1930.2307
$ws.Range("J63").Value = 7355.4
$ws.Range("K63").Value = 1930.2307
$ws.Range("L63").Value = 7355.4
$ws.Range("M63").Value = -1244.2307
$ws.Range("N63").Value = -8727.4
$ws.Range("H66").Value = 3437.2222
$ws.Range("I66").Value = 1930.2307
$ws.Range("J66").Value = 7355.4
$ws.Range("K66").Value = 9651.1535
$ws.Range("L66").Value = 36777
$ws.Range("M66").Value = -6219.1535
$ws.Range("N66").Value = -43641
$ws.Range("H121").Value = 0
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("H122").Value = 3599.1072
$ws.Range("I122").Value = 3604.32
$ws.Range("J122").Value = 3555.6667
$ws.Range("K122").Value = 10812.96
$ws.Range("L122").Value = 10667.0001
$ws.Range("M122").Value = -8362.960000000001
$ws.Range("N122").Value = -15567.0001
$ws.Range("H123").Value = 0
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("H124").Value = 16100.272
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 16100.272
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 16100.272
$ws.Range("N124").Value = -25920.272
$ws.Range("H125").Value = 59732.5
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 59732.5
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 59732.5
$ws.Range("N125").Value = -69572.5
$ws.Range("H126").Value = 1000000
$ws.Range("I126").Value = 1000000
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 3000000
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -2997530
$ws.Range("H127").Value = 0
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("H128").Value = 90000
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 90000
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 90000
$ws.Range("N128").Value = -99960
$ws.Range("H129").Value = 49999
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 49999
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 49999
$ws.Range("N129").Value = -59999
$ws.Range("H130").Value = 56594
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 56594
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 56594
$ws.Range("N130").Value = -66634
$ws.Range("H131").Value = 60469.875
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 60469.875
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 60469.875
$ws.Range("N131").Value = -70549.875
$ws.Range("H132").Value = 2623.1667
$ws.Range("I132").Value = 2139.5881
$ws.Range("J132").Value = 3255.5386
$ws.Range("K132").Value = 6418.7643
$ws.Range("L132").Value = 9766.6158
$ws.Range("M132").Value = -3888.7643
$ws.Range("N132").Value = -14826.6158
$ws.Range("H133").Value = 75065.25
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 75065.25
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 75065.25
$ws.Range("N133").Value = -80125.25
$ws.Range("H134").Value = 0
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("H135").Value = 44976.332
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 44976.332
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 44976.332
$ws.Range("N135").Value = -55116.332
$ws.Range("H137").Value = 100001
$ws.Range("I137").Value = 100001
$ws.Range("J137").Value = 100001
$ws.Range("K137").Value = 100001
$ws.Range("L137").Value = 100001
$ws.Range("M137").Value = -94901
$ws.Range("N137").Value = -110201
$ws.Range("H138").Value = 60825.715
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 60825.715
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 60825.715
$ws.Range("N138").Value = -71105.715
$ws.Range("H139").Value = 64500
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 64500
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 64500
$ws.Range("N139").Value = -74780
$ws.Range("H140").Value = 88538.164
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 88538.164
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 88538.164
$ws.Range("N140").Value = -98898.164
$ws.Range("H141").Value = 60235.266
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 60235.266
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 60235.266
$ws.Range("N141").Value = -70595.266

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3584.4146
$ws.Range("I86").Value = 1697.1578
$ws.Range("J86").Value = 5214.3184
$ws.Range("K86").Value = 1697.1578
$ws.Range("L86").Value = 5214.3184
$ws.Range("M86").Value = -574.1578
$ws.Range("N86").Value = -7460.3184
$ws.Range("H89").Value = 3584.4146
$ws.Range("I89").Value = 1697.1578
$ws.Range("J89").Value = 5214.3184
$ws.Range("K89").Value = 8485.789000000001
$ws.Range("L89").Value = 26071.592
$ws.Range("M89").Value = -2869.789000000001
$ws.Range("N89").Value = -37303.592
$ws.Range("H105").Value = 1689.8966
$ws.Range("I105").Value = 1688.7307
$ws.Range("K105").Value = 1688.7307
$ws.Range("M105").Value = 58.26929999999993
$ws.Range("H138").Value = 45776.668
$ws.Range("J138").Value = 45776.668
$ws.Range("L138").Value = 45776.668
$ws.Range("N138").Value = -56056.668
$ws.Range("H140").Value = 78860
$ws.Range("J140").Value = 78860
$ws.Range("L140").Value = 78860
$ws.Range("N140").Value = -89220

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 3894.4443
$ws.Range("I62").Value = 3881.25
$ws.Range("J62").Value = 4000
$ws.Range("K62").Value = 3881.25
$ws.Range("L62").Value = 4000
$ws.Range("M62").Value = -3257.25
$ws.Range("N62").Value = -5248
$ws.Range("H65").Value = 3894.4443
$ws.Range("I65").Value = 3881.25
$ws.Range("J65").Value = 4000
$ws.Range("K65").Value = 19406.25
$ws.Range("L65").Value = 20000
$ws.Range("M65").Value = -16286.25
$ws.Range("N65").Value = -26240

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 11974482
$ws.Range("I70").Value = 18150794
$ws.Range("J70").Value = 7877.25
$ws.Range("K70").Value = 18150794
$ws.Range("L70").Value = 7877.25
$ws.Range("M70").Value = -18150524
$ws.Range("N70").Value = -8417.25
$ws.Range("H73").Value = 11974482
$ws.Range("I73").Value = 18150794
$ws.Range("J73").Value = 7877.25
$ws.Range("K73").Value = 18150794
$ws.Range("L73").Value = 7877.25
$ws.Range("M73").Value = -18149858
$ws.Range("N73").Value = -9749.25
$ws.Range("H80").Value = 3212.5
$ws.Range("I80").Value = 3000
$ws.Range("J80").Value = 3340
$ws.Range("K80").Value = 3000
$ws.Range("L80").Value = 3340
$ws.Range("M80").Value = -2002
$ws.Range("N80").Value = -5336
$ws.Range("H83").Value = 3212.5
$ws.Range("I83").Value = 3000
$ws.Range("J83").Value = 3340
$ws.Range("K83").Value = 15000
$ws.Range("L83").Value = 16700
$ws.Range("M83").Value = -10008
$ws.Range("N83").Value = -26684

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1560.76
$ws.Range("I82").Value = 1643
$ws.Range("J82").Value = 1456.091
$ws.Range("K82").Value = 1643
$ws.Range("L82").Value = 1456.091
$ws.Range("M82").Value = -1282
$ws.Range("N82").Value = -2178.091
$ws.Range("H85").Value = 1560.76
$ws.Range("I85").Value = 1643
$ws.Range("J85").Value = 1456.091
$ws.Range("K85").Value = 1643
$ws.Range("L85").Value = 1456.091
$ws.Range("M85").Value = -395
$ws.Range("N85").Value = -3952.091

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H138").Value = 49525
$ws.Range("J138").Value = 49525
$ws.Range("L138").Value = 49525
$ws.Range("N138").Value = -59805
